$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=12; I='%'; J='Uninterpretable'},
    @{Row=15; I='aa'; J='Agree/Accept'},
    @{Row=21; I='sv'; J='Statement-opinion'},
    @{Row=24; I='sd'; J='Statement-non-opinion'},
    @{Row=28; I='aa'; J='Agree/Accept'},
    @{Row=32; I='aa'; J='Agree/Accept'},
    @{Row=43; I='aa'; J='Agree/Accept'},
    @{Row=66; I='sv'; J='Statement-opinion'},
    @{Row=73; I='sv'; J='Statement-opinion'},
    @{Row=75; I='sd'; J='Statement-non-opinion'},
    @{Row=97; I='ba'; J='Appreciation'},
    @{Row=109; I='sd'; J='Statement-non-opinion'},
    @{Row=112; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=115; I='aa'; J='Agree/Accept'},
    @{Row=116; I='sd'; J='Statement-non-opinion'},
    @{Row=119; I='sv'; J='Statement-opinion'},
    @{Row=123; I='sd'; J='Statement-non-opinion'},
    @{Row=124; I='sd'; J='Statement-non-opinion'},
    @{Row=129; I='%'; J='Uninterpretable'},
    @{Row=130; I='sd'; J='Statement-non-opinion'},
    @{Row=133; I='%'; J='Uninterpretable'},
    @{Row=136; I='sd'; J='Statement-non-opinion'},
    @{Row=137; I='sv'; J='Statement-opinion'},
    @{Row=139; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=145; I='sv'; J='Statement-opinion'},
    @{Row=153; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=156; I='sd'; J='Statement-non-opinion'},
    @{Row=162; I='sv'; J='Statement-opinion'},
    @{Row=164; I='%'; J='Uninterpretable'},
    @{Row=168; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=169; I='sd'; J='Statement-non-opinion'},
    @{Row=174; I='sv'; J='Statement-opinion'},
    @{Row=187; I='ba'; J='Appreciation'},
    @{Row=208; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=212; I='%'; J='Uninterpretable'},
    @{Row=214; I='%'; J='Uninterpretable'},
    @{Row=216; I='sd'; J='Statement-non-opinion'},
    @{Row=218; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=221; I='sd'; J='Statement-non-opinion'},
    @{Row=228; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=233; I='sd'; J='Statement-non-opinion'},
    @{Row=234; I='sd'; J='Statement-non-opinion'},
    @{Row=241; I='%'; J='Uninterpretable'},
    @{Row=244; I='sv'; J='Statement-opinion'},
    @{Row=246; I='sd'; J='Statement-non-opinion'},
    @{Row=251; I='sv'; J='Statement-opinion'},
    @{Row=270; I='sv'; J='Statement-opinion'},
    @{Row=285; I='aa'; J='Agree/Accept'},
    @{Row=293; I='sv'; J='Statement-opinion'},
    @{Row=305; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=318; I='sv'; J='Statement-opinion'},
    @{Row=323; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=328; I='sd'; J='Statement-non-opinion'},
    @{Row=329; I='sd'; J='Statement-non-opinion'},
    @{Row=353; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=356; I='ba'; J='Appreciation'},
    @{Row=371; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=381; I='sv'; J='Statement-opinion'},
    @{Row=387; I='ba'; J='Appreciation'},
    @{Row=393; I='sd'; J='Statement-non-opinion'},
    @{Row=396; I='sv'; J='Statement-opinion'},
    @{Row=404; I='sd'; J='Statement-non-opinion'},
    @{Row=420; I='sd'; J='Statement-non-opinion'},
    @{Row=424; I='ba'; J='Appreciation'},
    @{Row=427; I='aa'; J='Agree/Accept'},
    @{Row=445; I='%'; J='Uninterpretable'},
    @{Row=454; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=461; I='sd'; J='Statement-non-opinion'},
    @{Row=463; I='aa'; J='Agree/Accept'},
    @{Row=465; I='aa'; J='Agree/Accept'},
    @{Row=474; I='sd'; J='Statement-non-opinion'},
    @{Row=478; I='sd'; J='Statement-non-opinion'},
    @{Row=503; I='aa'; J='Agree/Accept'},
    @{Row=504; I='sd'; J='Statement-non-opinion'},
    @{Row=506; I='sd'; J='Statement-non-opinion'},
    @{Row=519; I='sd'; J='Statement-non-opinion'},
    @{Row=521; I='sd'; J='Statement-non-opinion'},
    @{Row=526; I='%'; J='Uninterpretable'},
    @{Row=552; I='sd'; J='Statement-non-opinion'},
    @{Row=553; I='sd'; J='Statement-non-opinion'},
    @{Row=560; I='%'; J='Uninterpretable'},
    @{Row=568; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=576; I='sv'; J='Statement-opinion'},
    @{Row=584; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=603; I='aa'; J='Agree/Accept'},
    @{Row=607; I='sd'; J='Statement-non-opinion'},
    @{Row=609; I='sd'; J='Statement-non-opinion'},
    @{Row=613; I='sv'; J='Statement-opinion'},
    @{Row=631; I='aa'; J='Agree/Accept'},
    @{Row=644; I='sd'; J='Statement-non-opinion'},
    @{Row=645; I='aa'; J='Agree/Accept'},
    @{Row=654; I='aa'; J='Agree/Accept'},
    @{Row=662; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=668; I='sd'; J='Statement-non-opinion'},
    @{Row=669; I='aa'; J='Agree/Accept'},
    @{Row=682; I='sd'; J='Statement-non-opinion'},
    @{Row=688; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Host "Applied $($updates.Count) updates"